$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.058095216751099
$ws.Range("B1").Value = 6.354019165039062
$ws.Range("C1").Value = 6.599575996398926
$ws.Range("D1").Value = 7.064857482910156
$ws.Range("E1").Value = 5.016282081604004
